# Update scripts with new TPM values (Icosl-Cd28, YoungD0 LR pairs)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (ECs -> FAPs target) ---
$ws.Range("D2").Value = "FAPs"
$ws.Range("G2").Value = 0.5407596666666666
$ws.Range("H2").Value = 1.622279
$ws.Range("I2").Value = 0.03618231591230665
$ws.Range("J2").Value = 0.03618231591230665
$ws.Range("M2").Value = 0.042039
$ws.Range("N2").Value = 0.126117
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.022732995627
$ws.Range("R2").Value = 0.204596960643
$ws.Range("S2").Value = 0.03618231591230665
$ws.Range("T2").Value = 0.03618231591230665

# --- Row 3 (FAPs -> FAPs, values refreshed) ---
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 11.32416633333333
$ws.Range("H3").Value = 33.972499
$ws.Range("I3").Value = 0.7577017832003755
$ws.Range("J3").Value = 0.7577017832003754
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.476056628487
$ws.Range("R3").Value = 4.284509656383
$ws.Range("S3").Value = 0.7577017832003755
$ws.Range("T3").Value = 0.7577017832003754

# --- Row 4 (FAPs -> MuSCs / ECs -> FAPs, values refreshed) ---
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 3.080487333333333
$ws.Range("H4").Value = 9.241461999999999
$ws.Range("I4").Value = 0.2061159008873179
$ws.Range("J4").Value = 0.2061159008873179
$ws.Range("M4").Value = 0.042039
$ws.Range("N4").Value = 0.126117
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.129500607006
$ws.Range("R4").Value = 1.165505463054
$ws.Range("S4").Value = 0.2061159008873179
$ws.Range("T4").Value = 0.2061159008873179

# --- Remove old rows 5-7 (table now only has rows 2-4 of data) ---
$ws.Range("A5:T7").Delete()
